$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.987.37"
$ws.Range("E2").Value = "  -1.33%  "
$ws.Range("D3").Value = "1.780.67"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("D4").Value = "'0.9996"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'316.00"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "'0.9997"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "'0.5392"
$ws.Range("E7").Value = "  -1.80%  "
$ws.Range("D8").Value = "'0.3772"
$ws.Range("E8").Value = "  -2.20%  "
$ws.Range("D9").Value = "'0.07439"
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("D10").Value = "'41.70"
$ws.Range("E10").Value = "  -2.10%  "
$ws.Range("D11").Value = "'1.093"
$ws.Range("E11").Value = "  -2.37%  "
$ws.Range("D12").Value = "'0.9998"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "'20.47"
$ws.Range("E13").Value = "  -3.19%  "
$ws.Range("D14").Value = "'6.078"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").Value = "'7.210"
$ws.Range("E15").Value = "  -1.90%  "
$ws.Range("D16").Value = "1.774.27"
$ws.Range("E16").Value = "  -1.52%  "
$ws.Range("D17").Value = "'88.28"
$ws.Range("E17").Value = "  -4.21%  "
$ws.Range("D18").Value = "'0.00001052"
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("D19").Value = "'0.06439"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "'0.9998"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "'17.22"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").Value = "'5.876"
$ws.Range("E22").Value = "  -1.92%  "
$ws.Range("D23").Value = "28.009.96"
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("E24").Value = "  -2.74%  "
$ws.Range("D25").Value = "'2.089"
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("D26").Value = "'156.01"
$ws.Range("E26").Value = "  -1.30%  "
$ws.Range("D27").Value = "'20.22"
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("D28").Value = "1.979.78"
$ws.Range("E28").Value = "  -1.47%  "
$ws.Range("D29").Value = "'2.290"
$ws.Range("E29").Value = "  -4.17%  "
$ws.Range("D30").Value = "'119.86"
$ws.Range("E30").Value = "  -3.16%  "
$ws.Range("E31").Value = "  -1.78%  "
$ws.Range("D32").Value = "'0.1049"
$ws.Range("E32").Value = "  +2.98%  "
$ws.Range("D33").Value = "'3.639"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "'5.511"
$ws.Range("E34").Value = "  -4.12%  "
$ws.Range("D35").Value = "'0.2256"
$ws.Range("E35").Value = "  -3.29%  "
$ws.Range("D36").Value = "'0.06423"
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("D37").Value = "'0.02272"
$ws.Range("E37").Value = "  -2.06%  "
$ws.Range("D38").Value = "'4.976"
$ws.Range("E38").Value = "  -1.77%  "
$ws.Range("D39").Value = "'8.431"
$ws.Range("E39").Value = "  -4.69%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "'11.12"
$ws.Range("E40").Value = "  -4.31%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6138"
$ws.Range("E41").Value = "  -4.22%  "
$ws.Range("D42").Value = "'1.430"
$ws.Range("E42").Value = "  +3.42%  "
$ws.Range("E43").Value = "  +1.62%  "
$ws.Range("D44").Value = "'0.9998"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "'13.23"
$ws.Range("E45").Value = "  -1.92%  "
$ws.Range("D46").Value = "'3.663"
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("D47").Value = "'0.5741"
$ws.Range("E47").Value = "  -3.93%  "
$ws.Range("D48").Value = "'126.26"
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("D49").Value = "'1.922"
$ws.Range("E49").Value = "  -3.17%  "
$ws.Range("D50").Value = "'1.180"
$ws.Range("E50").Value = "  +2.76%  "
$ws.Range("D51").Value = "'0.06788"
$ws.Range("E51").Value = "  -1.74%  "
